$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: extend the original chocolate-box vignette and add three new ---
# --- adaptation variants (columns C, D, E) as new scenario conditions.     ---

$ws.Range("B24").Value = "your friend has a box of chocolate. She knows that the round ones are filled with cream. She hands you a piece."
$ws.Range("C24").Value = "your friend has a box of with chocolates. She knows that 8 or 9 of the 10 have cream in them but she can’t tell from looking at them which ones are which. She hands you a piece."
$ws.Range("D24").Value = "your friend has a box of with chocolates. She knows that 5 or 6 of the 10 have cream in them but she can’t tell from looking at them which ones are which. She hands you a piece."
$ws.Range("E24").Value = "your friend has a box of with chocolates. She thinks that one of them is filled with cream but she can’t tell from looking at them which one it is. She hands you a piece."

# Row grows taller now that it wraps much longer text in every column.
$ws.Rows.Item(24).RowHeight = 64

# The view now scrolls one row further and the active selection moves to A24.
$ws.Range("A24").Select()
